# Updates crypto price/volume figures to the Jan 14 2023 12:52 UTC GitHub Actions refresh.
# Column D = Price, Column E = Volume(1h); both are plain-text cells (numeric-looking
# strings), so each new value is written with a leading apostrophe to force Excel to
# keep storing it as text instead of auto-converting it to a Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'303.64"
$ws.Range("E2").Value = "'5.83%"

# Row 3 - OKB
$ws.Range("D3").Value = "'31.99"
$ws.Range("E3").Value = "'8.67%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.246"
$ws.Range("E4").Value = "'3.43%"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.07514"
$ws.Range("E5").Value = "'11.74%"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "'7.849"
$ws.Range("E6").Value = "'7.09%"

# Row 7 - GateToken
$ws.Range("D7").Value = "'3.748"
$ws.Range("E7").Value = "'8.58%"

# Row 8 - FTXToken
$ws.Range("D8").Value = "'1.477"
$ws.Range("E8").Value = "'7.51%"

# Row 9 - MXToken
$ws.Range("D9").Value = "'0.9130"
$ws.Range("E9").Value = "'1.14%"

# Row 10 - One
$ws.Range("D10").Value = "'0.01671"
$ws.Range("E10").Value = "'2,479.13%"

# Row 11 - WazirX
$ws.Range("E11").Value = "'5.90%"

# Row 12 - LiechtensteinCryptoassetsExchange
$ws.Range("D12").Value = "'0.07436"
$ws.Range("E12").Value = "'7.34%"

# Row 13 - MandalaExchangeToken
$ws.Range("D13").Value = "'0.08002"
$ws.Range("E13").Value = "'4.88%"

# Row 14 - BitrueCoin
$ws.Range("D14").Value = "'0.03046"
$ws.Range("E14").Value = "'4.04%"

# Row 15 - BitMartToken
$ws.Range("D15").Value = "'0.09905"
$ws.Range("E15").Value = "'10.20%"

# Row 16 - BitForexToken
$ws.Range("D16").Value = "'0.001491"
$ws.Range("E16").Value = "'-4.88%"

# Row 17 - CoinExToken
$ws.Range("E17").Value = "'1.37%"

# Row 18 - TigerCash
$ws.Range("D18").Value = "'0.006180"
$ws.Range("E18").Value = "'-0.13%"

# Row 19 - LEO
$ws.Range("D19").Value = "'3.497"
$ws.Range("E19").Value = "'1.30%"

# Row 20 - BTSEToken
$ws.Range("D20").Value = "'2.237"
$ws.Range("E20").Value = "'0.29%"

# Row 21 - BitpandaEcosystemToken
$ws.Range("D21").Value = "'0.3310"
$ws.Range("E21").Value = "'3.29%"

# Row 22 - ProBitToken
$ws.Range("D22").Value = "'0.1339"
$ws.Range("E22").Value = "'1.84%"

# Row 23 - MCDex
$ws.Range("D23").Value = "'4.482"
$ws.Range("E23").Value = "'13.60%"

# Row 24 - ZBToken
$ws.Range("D24").Value = "'0.1624"
$ws.Range("E24").Value = "'4.07%"

# Row 25 - BitKan
$ws.Range("D25").Value = "'0.001212"
$ws.Range("E25").Value = "'0.86%"

# Row 26 - HotbitToken
$ws.Range("D26").Value = "'0.004451"
$ws.Range("E26").Value = "'1.93%"

# Row 27 - NitroEx
$ws.Range("D27").Value = "'0.0001396"
$ws.Range("E27").Value = "'16.29%"

# Row 28 - UpBots
$ws.Range("D28").Value = "'0.0001732"
$ws.Range("E28").Value = "'7.02%"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.04507"
$ws.Range("E40").Value = "'6.13%"

# Row 41 - KickToken
$ws.Range("D41").Value = "'0.007219"
$ws.Range("E41").Value = "'6.14%"

# Row 42 - BKEXToken
$ws.Range("D42").Value = "'0.1349"
$ws.Range("E42").Value = "'8.92%"

# Row 43 - CEJI
$ws.Range("D43").Value = "'0.002234"
$ws.Range("E43").Value = "'0.13%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.01397"

# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00006185"
$ws.Range("E45").Value = "'7.56%"

# Row 46 - BOLO
$ws.Range("E46").Value = "'-63.98%"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "'0.01294"
$ws.Range("E47").Value = "'-13.81%"
